# B6-PowerPoint.pptx edit:
#  1) Re-style the three data tables (slides 14-16) from the deck's
#     custom "Table_0" style to the built-in table style
#     {26EA0B1A-58A0-4A0E-B2AB-55D307264C82}.
#  2) Swap the presentation's theme palette: the deck's live theme
#     (Integral / "Red Violet") is replaced by the stock Office theme
#     palette (the two <a:theme> parts effectively trade places).

$p = $ppt.ActivePresentation

# --- 1. Table styles -------------------------------------------------
$newStyle = "{26EA0B1A-58A0-4A0E-B2AB-55D307264C82}"
foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newStyle)
        }
    }
}

# --- 2. Theme colors ---------------------------------------------------
# Index order matches MsoThemeColorSchemeIndex:
#  1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#  8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$master = $p.SlideMaster
$tcs = $master.Theme.ThemeColorScheme

$tcs.Item(1).RGB  = 0          # dk1      000000
$tcs.Item(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388    # dk2      44546A
$tcs.Item(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  4472C4
$tcs.Item(10).RGB = 4697456    # accent6  70AD47
$tcs.Item(11).RGB = 12673797   # hlink    0563C1
$tcs.Item(12).RGB = 7491477    # folHlink 954F72
